$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7: fill in Increase/Decrease column (I7) -> "increase " (trailing space)
$ws.Range("I7").Value = "increase "

# Row 6: fill in Increase/Decrease column (I6) -> "increase"
$ws.Range("I6").Value = "increase"

# Row 8: fill in Increase/Decrease column (I8) -> "increase"
$ws.Range("I8").Value = "increase"

# Row 9: update values for a test case that decreases
$ws.Range("C9").Value = 2
$ws.Range("G9").Value = 57878572
$ws.Range("H9").Value = 275221788
$ws.Range("I9").Value = "decrease"

# Row 10: fill in Population Change (G10) and Increase/Decrease (I10)
$ws.Range("G10").Value = 1088743
$ws.Range("I10").Value = "decrease"

# Update selection to match final state
$ws.Range("H13").Select()
